# Fruta / hortaliza, semanal
# A new weekly price-report row needs to be inserted into the daily
# logic table for "Femacal de La Calera - Tuna". Existing rows 29..111
# shift down to 30..112, and the newly opened row 29 is populated with
# the latest report's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 29..111 down to 30..112, opening up a blank row 29
# (carries over the D-column date style automatically).
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record.
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44623
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107011
$ws.Range("J29").Value = "Tuna"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 70
$ws.Range("N29").Value = 14000
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 14000
$ws.Range("Q29").Value = "$/caja 16 kilos"
$ws.Range("R29").Value = "Cabildo"
$ws.Range("S29").Value = 875
$ws.Range("T29").Value = 16
